# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# F2: 334 -> 335
# F3: 1362 -> 1367

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 335
    $ws.Range("F3").Value = 1367
}
